$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "O" column labels (columns O1:O12) to the new, simplified/lowercased
# label names, in the same row order as before.
$newLabels = @(
    "num_employess",
    "it_job",
    "employer_awareness",
    "family_history",
    "past_history",
    "current_state",
    "treatment",
    "age",
    "gender",
    "residence",
    "work_location",
    "remote"
)

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = $i + 1
    $ws.Range("O$row").Value = $newLabels[$i]
}

# Remove the now-unused columns between the question text (column A) and the
# label column, so the label column ends up as column F.
$ws.Columns("B:J").Delete()

# Make column A wide enough to show the full question text (matches the
# "best fit" width Excel computed for the longest question string).
$ws.Columns("A:A").ColumnWidth = 102

# Update the active selection.
$ws.Range("A14").Select() | Out-Null
